$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-DateText($cellRef, $text) {
    # A couple of these day/month strings are ambiguous (day <= 12) and
    # would otherwise get auto-parsed into a real date serial by the
    # engine's locale-aware literal-entry heuristics. The leading
    # apostrophe forces literal text entry (same as typing it into Excel),
    # so the stored value stays exactly the original "dd-mm-yyyy" string.
    # ClearFormats() afterwards drops the transient "text-entry" style bit
    # so the cell ends up with no explicit style, same as before the edit.
    $rng = $ws.Range($cellRef)
    $rng.Value = "'" + $text
    $rng.ClearFormats()
}

# Row 3: date format change + value updates
Set-DateText "A3" "28-07-2022"
$ws.Range("D3").Value = 1
$ws.Range("G3").Value = 1

# Row 4: date format change only
Set-DateText "A4" "01-08-2022"

# Row 5: date format change only
Set-DateText "A5" "04-08-2022"

# Row 6: date format change only
Set-DateText "A6" "08-08-2022"

# Row 7: date format change only
Set-DateText "A7" "11-08-2022"

# Row 8: date format change only
Set-DateText "A8" "15-08-2022"

# Row 9: date format change only
Set-DateText "A9" "18-08-2022"

# Row 10: date format change + value updates
Set-DateText "A10" "22-08-2022"
$ws.Range("D10").Value = 1
$ws.Range("G10").Value = 1

# Row 11: date format change only
Set-DateText "A11" "25-08-2022"

# Row 12: date format change + value updates
Set-DateText "A12" "29-08-2022"
$ws.Range("D12").Value = 1
$ws.Range("E12").Value = 1
$ws.Range("H12").Value = 0

# Row 13: date format change + value updates
Set-DateText "A13" "01-09-2022"
$ws.Range("D13").Value = 1
$ws.Range("E13").Value = 1
$ws.Range("H13").Value = 0

# Row 14: date format change only
Set-DateText "A14" "05-09-2022"

# Row 15: date format change only
Set-DateText "A15" "08-09-2022"

# Row 16: date format change only
Set-DateText "A16" "12-09-2022"

# Row 17: date format change only
Set-DateText "A17" "15-09-2022"

# Row 18: date format change only
Set-DateText "A18" "19-09-2022"

# Row 19: date format change only
Set-DateText "A19" "22-09-2022"

# Row 20: date format change only
Set-DateText "A20" "26-09-2022"

# Row 21: date format change only
Set-DateText "A21" "29-09-2022"
